$wb = $excel.ActiveWorkbook

# Rename the first sheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Variables"

# Add a new worksheet after the first one, named "Categories"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Categories"

# Apply the same header style used on Variables!A1:C1 (bold, centered) to the
# new sheet's header row, then fill in the header labels.
$ws1.Range("A1:C1").Copy() | Out-Null
$ws2.Range("A1:C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws2.Range("A1").Value = "variable"
$ws2.Range("B1").Value = "name"
$ws2.Range("C1").Value = "label"

# Update selection on Categories sheet to E12 first (it becomes active on Add)
$ws2.Range("E12").Select() | Out-Null

# Update selection on Variables sheet to B30 and make it the active sheet
$ws1.Select() | Out-Null
$ws1.Range("B30").Select() | Out-Null

$wb.Save() | Out-Null
